$newRowValues = @{
    3 = @{ "A"="n"; "B"="Barbara Bee"; "C"="Barbara"; "D"="Bee"; "E"="Mrs"; "F"="Barbara.Bee@hotmail.com"; "G"="Password123"; "H"="Password125"; "I"="Barbara.Beehotmail.com"; "J"="XYZ"; "K"="574-897-2346"; "L"="124 Main St"; "M"="Clain"; "N"="VA"; "O"="22315"; "P"="United States"; "Q"="5"; "R"="May"; "S"="1987"; "T"="2test@hotmail.com"; "U"="We Love Java" }
    4 = @{ "A"="n"; "B"="Mostafa King"; "C"="Mostafa"; "D"="King"; "E"="Mr"; "F"="Mostafa.King2@hotmail.com"; "G"="Password123"; "H"="Password128"; "I"="Mostafa.Kinghotmail.com"; "J"="XYZ"; "K"="574-897-2349"; "L"="127 Main St"; "M"="Clain"; "N"="VA"; "O"="22318"; "P"="India"; "Q"="10"; "R"="August"; "S"="1989"; "T"="5test@hotmail.com"; "U"="We Love Java" }
    5 = @{ "A"="n"; "B"="Kevin Lee"; "C"="Kevin"; "D"="Lee"; "E"="Mr"; "F"="Kevin.Lee7@hotmail.com"; "G"="Password123"; "H"="Password124"; "I"="Kevin.Lee10hotmail.com"; "J"="XYZ"; "K"="574-897-2345"; "L"="123 Main St"; "M"="Clain"; "N"="VA"; "O"="22314"; "P"="Canada"; "Q"="3"; "R"="January"; "S"="2001"; "T"="1test@hotmail.com"; "U"="We Love Java" }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customer")
$ws.Activate()

foreach ($r in $newRowValues.Keys) {
    $rowHash = $newRowValues[$r]
    foreach ($c in $rowHash.Keys) {
        $ws.Range("$c$r").Value = $rowHash[$c]
    }
}

# Update the active selection on the "customer" sheet
$ws.Range("F12").Select()

# Update the saved window position/size for the workbook
$win = $excel.ActiveWindow
$win.Left = 1497.75
$win.Top = 109.5
$win.Width = 1080
$win.Height = 569.25
